# Update gh-pages output data (generated at 456a3b4).
# Applies the same set of "想去人数" / "最低票价" refreshes to both the
# "展览" sheet and the "全部类型" sheet (the latter repeats the same rows
# at a two-row offset because it also contains rows from other sheets).

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet {
    param($ws, $rowG3, $rowF4, $rowG5, $rowF7, $rowF8, $rowF9, $rowF10, $rowF13, $rowF14, $rowF17, $rowF18, $rowF19, $rowF20, $rowF21, $rowF23)

    $ws.Cells.Item($rowG3, 7).Value = 60

    $ws.Cells.Item($rowF4, 6).Value = 8539
    $ws.Cells.Item($rowF4, 7).Value = 75

    $ws.Cells.Item($rowG5, 7).Value = "不可售"

    $ws.Cells.Item($rowF7, 6).Value = 21
    $ws.Cells.Item($rowF8, 6).Value = 80
    $ws.Cells.Item($rowF9, 6).Value = 1331
    $ws.Cells.Item($rowF10, 6).Value = 106
    $ws.Cells.Item($rowF13, 6).Value = 9182
    $ws.Cells.Item($rowF14, 6).Value = 150
    $ws.Cells.Item($rowF17, 6).Value = 169
    $ws.Cells.Item($rowF18, 6).Value = 342
    $ws.Cells.Item($rowF19, 6).Value = 6078
    $ws.Cells.Item($rowF20, 6).Value = 1042
    $ws.Cells.Item($rowF21, 6).Value = 59
    $ws.Cells.Item($rowF23, 6).Value = 102
}

# "展览" sheet: rows line up 1:1 with the logical row numbers.
$wsExpo = $wb.Worksheets.Item("展览")
Update-ExpoSheet $wsExpo 3 4 5 7 8 9 10 13 14 17 18 19 20 21 23

# "全部类型" sheet: same events, but shifted down by 2 rows starting at row 13
# because two extra "演出" rows are interleaved above them.
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExpoSheet $wsAll 3 4 5 7 8 9 10 15 16 19 20 21 22 23 25
